$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties, matching the style of the
# existing header row (bold font, borders, centered alignment).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Team record columns for every data row (2-55): Wins=78, Losses=84, Ties=0
$ws.Range("AD2:AD55").Value = 78
$ws.Range("AE2:AE55").Value = 84
$ws.Range("AF2:AF55").Value = 0
